# Refresh cryptocurrency snapshot values (D = Price, E = Volume(1h) change)
# for every row whose source data changed in this run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param(
        [string]$Address,
        [string]$Text
    )
    $range = $ws.Range($Address)
    $looksNumeric = $Text -match "^[+-]?[0-9]*\.?[0-9]+$"
    if ($looksNumeric) {
        # Plain numeric-looking strings (e.g. "571.13") must stay text,
        # matching the sheet's existing inline-string "Price" column --
        # force text via NumberFormat, then restore the default
        # (unstyled) cell style so formatting is left unchanged.
        $range.NumberFormat = "@"
        $range.Value = $Text
        $range.Style = "Normal"
    } else {
        $range.Value = $Text
    }
}

$updates = @(
    @{ Address = "D2"; Text = "60.988.57" },
    @{ Address = "E2"; Text = "  +0.48%  " },
    @{ Address = "D3"; Text = "3.384.24" },
    @{ Address = "E3"; Text = "  +0.06%  " },
    @{ Address = "E4"; Text = "  -0.01%  " },
    @{ Address = "D5"; Text = "571.13" },
    @{ Address = "E5"; Text = "  +0.00%  " },
    @{ Address = "D6"; Text = "141.88" },
    @{ Address = "E6"; Text = "  +0.33%  " },
    @{ Address = "E7"; Text = "  -0.02%  " },
    @{ Address = "E8"; Text = "  +0.26%  " },
    @{ Address = "D9"; Text = "7.64" },
    @{ Address = "E9"; Text = "  +1.90%  " },
    @{ Address = "E10"; Text = "  -0.56%  " },
    @{ Address = "E11"; Text = "  -1.41%  " },
    @{ Address = "D12"; Text = "3.965.24" },
    @{ Address = "E12"; Text = "  +0.05%  " },
    @{ Address = "E13"; Text = "  +1.80%  " },
    @{ Address = "D14"; Text = "27.91" },
    @{ Address = "E14"; Text = "  -0.35%  " },
    @{ Address = "D15"; Text = "0.0000171" },
    @{ Address = "E15"; Text = "  +0.55%  " },
    @{ Address = "D16"; Text = "3.390.51" },
    @{ Address = "E16"; Text = "  +0.04%  " },
    @{ Address = "D17"; Text = "61.089.34" },
    @{ Address = "E17"; Text = "  +0.41%  " },
    @{ Address = "D18"; Text = "6.08" },
    @{ Address = "E18"; Text = "  -3.13%  " },
    @{ Address = "D19"; Text = "13.60" },
    @{ Address = "E19"; Text = "  -3.47%  " },
    @{ Address = "D20"; Text = "8.89" },
    @{ Address = "E20"; Text = "  -0.54%  " },
    @{ Address = "D21"; Text = "383.99" },
    @{ Address = "E21"; Text = "  -1.22%  " },
    @{ Address = "D22"; Text = "75.19" },
    @{ Address = "E22"; Text = "  +2.64%  " },
    @{ Address = "E23"; Text = "  -1.34%  " },
    @{ Address = "E24"; Text = "  +0.00%  " },
    @{ Address = "E25"; Text = "  -1.57%  " },
    @{ Address = "D26"; Text = "3.524.57" },
    @{ Address = "E26"; Text = "  +0.01%  " },
    @{ Address = "E27"; Text = "  +2.49%  " },
    @{ Address = "E28"; Text = "  -0.03%  " },
    @{ Address = "D29"; Text = "7.26" },
    @{ Address = "E29"; Text = "  -1.49%  " },
    @{ Address = "D30"; Text = "7.95" },
    @{ Address = "E30"; Text = "  -1.49%  " },
    @{ Address = "E31"; Text = "  +0.23%  " },
    @{ Address = "E33"; Text = "  -4.47%  " },
    @{ Address = "D34"; Text = "23.20" },
    @{ Address = "E34"; Text = "  -2.41%  " },
    @{ Address = "D35"; Text = "6.93" },
    @{ Address = "E35"; Text = "  +0.29%  " },
    @{ Address = "D36"; Text = "166.43" },
    @{ Address = "E36"; Text = "  -0.22%  " },
    @{ Address = "D37"; Text = "3.418.30" },
    @{ Address = "E37"; Text = "  +0.25%  " },
    @{ Address = "D38"; Text = "4.97" },
    @{ Address = "E38"; Text = "  -1.30%  " },
    @{ Address = "E39"; Text = "  -3.01%  " },
    @{ Address = "D40"; Text = "0.0766" },
    @{ Address = "E40"; Text = "  -1.42%  " },
    @{ Address = "D41"; Text = "26.89" },
    @{ Address = "E41"; Text = "  -0.04%  " },
    @{ Address = "E42"; Text = "  -0.04%  " },
    @{ Address = "E43"; Text = "  -0.59%  " },
    @{ Address = "E44"; Text = "  -1.92%  " },
    @{ Address = "E45"; Text = "  -1.62%  " },
    @{ Address = "E46"; Text = "  +0.03%  " },
    @{ Address = "D47"; Text = "2.448.17" },
    @{ Address = "E47"; Text = "  -3.86%  " },
    @{ Address = "D48"; Text = "22.89" },
    @{ Address = "E48"; Text = "  +0.21%  " },
    @{ Address = "D49"; Text = "6.68" },
    @{ Address = "E49"; Text = "  -2.15%  " },
    @{ Address = "E50"; Text = "  +10.49%  " },
    @{ Address = "E51"; Text = "  +0.84%  " }
)

foreach ($u in $updates) {
    Set-TextCell $u.Address $u.Text
}
